# Generate Report for Archive
# - Flip the localization "Status" from "Ready for handoff" to "In Translation"
#   everywhere it appears (Overview!E2:F4 and the per-locale tables' Status
#   column, zh-cn!C2:C4 / de-de!C2:C4).
# - The Status column narrows to fit the new (shorter) text, so re-fit those
#   columns afterwards.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn status lives in column E, de-de status in column F.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus

# Per-locale detail sheets: Status lives in column C.
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("C4").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("C4").Value = $newStatus

# Re-fit the Status columns to the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
